$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 used to be "Pulley pitch diameter (um)" = B5*10000.
# It now becomes "Pulley circumference (m)" = B5*PI().
$ws.Range("A6").Value = "Pulley circumference (m)"
$ws.Range("B6").Formula = "=B5*PI()"

# Row 10 ("Max encoder ticks") now divides by the pulley circumference (B6)
# instead of the pulley pitch diameter (B5), to account for the change in direction.
$ws.Range("B10").Formula = "=B9/B6*B8"

# Update the selection to match the recorded view state.
$ws.Range("E9").Select()
